# Updated cryptos list — refreshes the "Price" (column D) and
# "Volume(1h)" (column E) figures for every coin row (rows 2-51) on
# Sheet1 with the latest scraped snapshot.
#
# Both columns are stored as literal text in the workbook (e.g. prices
# like "65.159.36" use '.' as a thousands separator rather than a
# decimal point, and volumes are padded percent strings such as
# "  +1.36%  "), so every write below forces the destination cell to
# Text format first and clears that formatting override again right
# after the write. That stops Excel's COM layer from "helpfully"
# re-interpreting values that happen to look numeric (e.g. "1.00",
# "0.637") as actual numbers, while leaving the cell's style untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '65.159.36' },
    @{ Cell = 'E2'; Value = '  +1.36%  ' },
    @{ Cell = 'D3'; Value = '3.382.04' },
    @{ Cell = 'E3'; Value = '  +1.08%  ' },
    @{ Cell = 'E4'; Value = '  -0.06%  ' },
    @{ Cell = 'D5'; Value = '555.75' },
    @{ Cell = 'E5'; Value = '  +0.10%  ' },
    @{ Cell = 'D6'; Value = '174.94' },
    @{ Cell = 'E6'; Value = '  -0.05%  ' },
    @{ Cell = 'D7'; Value = '0.632' },
    @{ Cell = 'E7'; Value = '  +2.24%  ' },
    @{ Cell = 'D8'; Value = '3.373.42' },
    @{ Cell = 'E8'; Value = '  +1.16%  ' },
    @{ Cell = 'E9'; Value = '  -0.11%  ' },
    @{ Cell = 'E10'; Value = '  +6.17%  ' },
    @{ Cell = 'D11'; Value = '0.637' },
    @{ Cell = 'E11'; Value = '  +1.55%  ' },
    @{ Cell = 'D12'; Value = '53.71' },
    @{ Cell = 'E12'; Value = '  -1.34%  ' },
    @{ Cell = 'E13'; Value = '  +2.55%  ' },
    @{ Cell = 'D14'; Value = '9.17' },
    @{ Cell = 'E14'; Value = '  +1.26%  ' },
    @{ Cell = 'D15'; Value = '3.923.21' },
    @{ Cell = 'E15'; Value = '  +0.62%  ' },
    @{ Cell = 'D16'; Value = '18.32' },
    @{ Cell = 'E16'; Value = '  -0.34%  ' },
    @{ Cell = 'D17'; Value = '3.402.10' },
    @{ Cell = 'E17'; Value = '  +1.22%  ' },
    @{ Cell = 'E18'; Value = '  +0.27%  ' },
    @{ Cell = 'D19'; Value = '65.191.78' },
    @{ Cell = 'E19'; Value = '  +1.33%  ' },
    @{ Cell = 'D20'; Value = '11.85' },
    @{ Cell = 'E20'; Value = '  +0.35%  ' },
    @{ Cell = 'E21'; Value = '  +1.76%  ' },
    @{ Cell = 'D22'; Value = '455.89' },
    @{ Cell = 'E22'; Value = '  +0.39%  ' },
    @{ Cell = 'D23'; Value = '4.89' },
    @{ Cell = 'E23'; Value = '  +0.30%  ' },
    @{ Cell = 'D24'; Value = '14.30' },
    @{ Cell = 'E24'; Value = '  +7.59%  ' },
    @{ Cell = 'E25'; Value = '  +0.27%  ' },
    @{ Cell = 'D26'; Value = '87.32' },
    @{ Cell = 'E26'; Value = '  +1.94%  ' },
    @{ Cell = 'D27'; Value = '2.89' },
    @{ Cell = 'E27'; Value = '  +1.45%  ' },
    @{ Cell = 'D28'; Value = '10.70' },
    @{ Cell = 'E28'; Value = '  -1.95%  ' },
    @{ Cell = 'D29'; Value = '8.72' },
    @{ Cell = 'E29'; Value = '  -0.33%  ' },
    @{ Cell = 'D30'; Value = '31.11' },
    @{ Cell = 'E30'; Value = '  +4.00%  ' },
    @{ Cell = 'D31'; Value = '6.54' },
    @{ Cell = 'E31'; Value = '  -1.13%  ' },
    @{ Cell = 'D32'; Value = '63.33' },
    @{ Cell = 'E32'; Value = '  +8.16%  ' },
    @{ Cell = 'E33'; Value = '  +0.14%  ' },
    @{ Cell = 'D34'; Value = '579.79' },
    @{ Cell = 'E34'; Value = '  -0.70%  ' },
    @{ Cell = 'E35'; Value = '  -0.09%  ' },
    @{ Cell = 'E36'; Value = '  +0.12%  ' },
    @{ Cell = 'D37'; Value = '3.63' },
    @{ Cell = 'E37'; Value = '  +4.00%  ' },
    @{ Cell = 'D38'; Value = '0.143' },
    @{ Cell = 'E38'; Value = '  +1.87%  ' },
    @{ Cell = 'D39'; Value = '35.73' },
    @{ Cell = 'E39'; Value = '  +0.14%  ' },
    @{ Cell = 'E40'; Value = '  -0.19%  ' },
    @{ Cell = 'D41'; Value = '0.0₃0739' },
    @{ Cell = 'E41'; Value = '  -1.93%  ' },
    @{ Cell = 'D42'; Value = '3.088.62' },
    @{ Cell = 'E42'; Value = '  -0.26%  ' },
    @{ Cell = 'E43'; Value = '  +1.85%  ' },
    @{ Cell = 'D44'; Value = '2.76' },
    @{ Cell = 'E44'; Value = '  -0.86%  ' },
    @{ Cell = 'E45'; Value = '  +2.86%  ' },
    @{ Cell = 'E46'; Value = '  -2.82%  ' },
    @{ Cell = 'D47'; Value = '3.17' },
    @{ Cell = 'E47'; Value = '  -1.64%  ' },
    @{ Cell = 'D48'; Value = '1.00' },
    @{ Cell = 'E48'; Value = '  -0.19%  ' },
    @{ Cell = 'D49'; Value = '142.29' },
    @{ Cell = 'E49'; Value = '  +4.99%  ' },
    @{ Cell = 'E50'; Value = '  -1.73%  ' },
    @{ Cell = 'E51'; Value = '  +0.15%  ' }

)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.ClearFormats()
}
